# Finish anti-viral therapy update (#176)
# - fix swapped chemical_identity / defined_class values on row 42
# - add anti-HSV agent / anti-HSV-1 agent / anti-HSV-2 agent / glucocorticoid rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix row 42: A42 and C42 had been swapped ---
$ws.Range("A42").Value = "MAXO_0000249"
$ws.Range("C42").Value = "CHEBI_59886"

# --- row 44: anti-HSV agent therapy (taller row / larger font on C44) ---
$ws.Range("A44").Value = "MAXO_0000637"
$ws.Range("B44").Value = "anti-HSV agent therapy"
$ws.Range("C44").Value = "CHEBI:64952"
$ws.Range("D44").Value = "anti-HSV agent"

# --- row 45: anti-HSV-1 agent therapy ---
$ws.Range("A45").Value = "MAXO_0000638"
$ws.Range("B45").Value = "anti-HSV-1 agent therapy"
$ws.Range("C45").Value = "CHEBI:64953"
$ws.Range("D45").Value = "anti-HSV-1 agent"

# --- row 46: anti-HSV-2 agent therapy ---
$ws.Range("A46").Value = "MAXO_0000639"
$ws.Range("B46").Value = "anti-HSV-2 agent therapy"
$ws.Range("C46").Value = "CHEBI:64954"
$ws.Range("D46").Value = "anti-HSV-2 agent"

# --- row 47: glucocorticoid agent therapy ---
$ws.Range("A47").Value = "MAXO_0000220"
$ws.Range("B47").Value = "glucocorticicoid agent therapy"
$ws.Range("C47").Value = "CHEBI:24261"
$ws.Range("D47").Value = "glucocorticoid"

# --- formatting: C45 / C46 reuse the existing "big Helvetica" style (already
#     used by e.g. C4), copied over without disturbing the pasted values ---
$ws.Range("C4").Copy()
$ws.Range("C45").PasteSpecial(-4122)
$ws.Range("C46").PasteSpecial(-4122)

# --- formatting: C44 needs a brand new font (sz 13 Helvetica); build it on a
#     scratch cell far outside the used range, copy the format across, then
#     remove the scratch row entirely ---
$helper = $ws.Range("A500")
$helper.Font.Name = "Helvetica"
$helper.Font.Size = 13
$helper.Copy()
$ws.Range("C44").PasteSpecial(-4122)
$helper.EntireRow.Delete()

# --- row heights for the new rows ---
$ws.Range("A44:D44").RowHeight = 17
$ws.Range("A45:D45").RowHeight = 21
$ws.Range("A46:D46").RowHeight = 21

# --- sheet view bookkeeping to match the saved workbook state ---
$ws.Range("B48").Select()
$excel.ActiveWindow.ScrollRow = 14

Write-Host "edit applied"
